# Commit: "Added abbr for Speed of germination"
#
# Adds mathematical symbol abbreviations to three index names on Sheet1
# (Speed of germination / Speed of accumulated germination / Corrected
# speed of germination rows) and appends the companion formula describing
# the "corrected speed of accumulated germination" to the Details cell.
# Also nudges the visible selection down towards those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 - "Speed of germination" index name gains the ($S$) symbol.
$ws.Range("A18").Value = @'
Speed of germination ($S$) or Germination rate Index or index of velocity of germination or Emergence rate index (Allan, Vogel and Peterson; Erbach; Hsu and Nelson) or Germination index (AOSA)
'@

# Row 19 - "Speed of accumulated germination" gains the ($S_{accumulated}$) symbol.
$ws.Range("A19").Value = @'
Speed of accumulated germination ($S_{accumulated}$)
'@

# Row 20 - renamed / gains the ($\hat{S}$) symbol.
$ws.Range("A20").Value = @'
Corrected speed of germination or Corrected germination rate index ($\hat{S}$)
'@

# Row 20 Details - formula explanation extended with the accumulated-germination variant.
$ws.Range("C20").Value = @'
It is computed as follows.
$$S_{corrected} = \frac{S}{FGP}$$
Where, $S$ is the germination speed computed with germination percentage instead of counts and $FGP$ is the final germination percentage or germinability.
It can also be computed from speed of accumulated germination (computed with germination percentage).
$$\hat{S}_{accumulated} = \frac{S_{accumulated}}{FGP}$$
Where, $S_{accumulated}$ is the speed of accumulated germination computed with germination percentage instead of counts and $FGP$ is the final germination percentage or germinability.
'@

# Row 20 now holds considerably more text, matching the taller row in the edit.
$ws.Rows.Item(20).RowHeight = 150

# The author's cursor ended up at C24 after editing, with the view scrolled
# so row 20 is the top visible row.
$ws.Range("C24").Select()
